$d = $word.ActiveDocument

# 1. Contact info line: remove phone number, then bump font size 18 -> 20
$d.Content.Find.Execute(
    "annacape@colostate.edu | 970.227.3390 | https://www.linkedin.com/in/anna-capels-204327276/ | https://github.com/AnnaC-1 | https://acapels.com/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "annacape@colostate.edu | https://www.linkedin.com/in/anna-capels-204327276/ | https://github.com/AnnaC-1 | https://acapels.com/",
    2)

# Resize that paragraph's text (and paragraph mark) from 9pt/18 half-points
# to 10pt/20 half-points (covers both w:sz and w:szCs).
$contactPara = $d.Paragraphs(3)
$paraRange = $contactPara.Range
$paraRange.Font.SizeBi = 10
$paraRange.Font.Size = 10

# 2. Merge "Python (Pandas, NumPy, Matplotlib.pyplot, Seaborn)"
$d.Content.Find.Execute(
    "Python (Pandas, NumPy, Matplotlib.pyplot, Seaborn)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Python (Pandas, NumPy, Matplotlib.pyplot, Seaborn)",
    2)

# 3. Merge "-Utilized Yfinance to perform..."
$d.Content.Find.Execute(
    "-Utilized Yfinance to perform stock risk analysis by data processing and visualization with financial techniques ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-Utilized Yfinance to perform stock risk analysis by data processing and visualization with financial techniques ",
    2)

# 4. Merge "VaR for easy communication of financial insights". A plain
# Find/Replace leaves a dangling <w:proofErr w:type="spellStart"/> behind
# (it sits right at the paragraph's start, ahead of the matched text, so a
# text-only replace never touches it). Rebuild the paragraph body via
# Range.InsertXML instead, which drops the stray proofErr markers cleanly
# while keeping the paragraph's indentation/formatting (pPr) intact.
$varPara = $d.Paragraphs(17)
$varRange = $varPara.Range
$varRange.End = $varRange.End - 1
$varRange.Delete()
$varXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="432" w:firstLine="288"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>VaR for easy communication of financial insights</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$varRange.InsertXML($varXml)

# 5. Merge "Python (TensorFlow, Keras)"
$d.Content.Find.Execute(
    "Python (TensorFlow, Keras)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Python (TensorFlow, Keras)",
    2)

# 6. Merge "- Spanish | Written, Reading, and Verbally Fluent"
$d.Content.Find.Execute(
    "- Spanish | Written, Reading, and Verbally Fluent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Spanish | Written, Reading, and Verbally Fluent",
    2)
